# Applies the "updated outcomes in iteration 1 plan" edit:
#  1) Heading 1 "System Interfaces" - swap explicit Times New Roman /
#     b=0 / bCs / sz=20 run-props for theme (minorHAnsi) fonts + szCs=24,
#     both on the paragraph mark and on the "System Interfaces" run.
#  2) "4.0.4 User Interfaces " paragraph - drop the gramStart/gramEnd
#     proofing-error bookmarks around the word and lower-case it
#     ("For" -> "for").
#  3) Heading 2 "Interfaces to External Systems or Devices" - same
#     Times New Roman -> theme-font swap as (1), minus the size tweak,
#     across the paragraph mark and both runs.

$d = $word.ActiveDocument

$pkgHead = '<?xml version="1.0" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" mc:Ignorable="w14 w15">' +
    '<w:body>'
$pkgTail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- 1) Heading 1 "System Interfaces" ------------------------------------

Write-Host "Updating 'System Interfaces' heading..."
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 1" -and $p.Range.Text.Trim() -eq "System Interfaces") {
        $xml = '<w:p w14:paraId="1B06DE97" w14:textId="77777777" w:rsidR="00E27454" w:rsidRPr="00E27454" w:rsidRDefault="008C2B65" w:rsidP="00E27454"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00C00447"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:br w:type="page"/></w:r><w:bookmarkStart w:id="2" w:name="_Toc492960765"/><w:r w:rsidR="00E27454" w:rsidRPr="00E27454"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>System Interfaces</w:t></w:r></w:p>'
        $p.Range.InsertXML($pkgHead + $xml + $pkgTail)
        Write-Host "  done."
        break
    }
}

# --- 2) "4.0.4 User Interfaces For ..." -> "for" -------------------------

Write-Host "Updating '4.0.4 User Interfaces' paragraph..."
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "4.0.4 User Interfaces*") {
        $xml = '<w:p w14:paraId="60407F34" w14:textId="77777777" w:rsidR="00E27454" w:rsidRPr="00E27454" w:rsidRDefault="00E27454" w:rsidP="00E27454"><w:pPr><w:rPr><w:bCs/><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:r w:rsidRPr="00E27454"><w:rPr><w:bCs/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve">4.0.4 User Interfaces </w:t></w:r><w:r w:rsidRPr="00E27454"><w:rPr><w:bCs/><w:lang w:val="en-AU"/></w:rPr><w:t>for</w:t></w:r><w:r w:rsidRPr="00E27454"><w:rPr><w:bCs/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> the ADHD management mobile application, user interfaces should prioritize simplicity, user-friendliness, and accessibility. Consider the following requirements:</w:t></w:r></w:p>'
        $p.Range.InsertXML($pkgHead + $xml + $pkgTail)
        Write-Host "  done."
        break
    }
}

# --- 3) Heading 2 "Interfaces to External Systems or Devices" ------------

Write-Host "Updating 'Interfaces to External Systems or Devices' heading..."
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 2" -and $p.Range.Text.Trim() -eq "Interfaces to External Systems or Devices") {
        $xml = '<w:p w14:paraId="76ABEE37" w14:textId="77777777" w:rsidR="00E27454" w:rsidRPr="00E27454" w:rsidRDefault="00E27454" w:rsidP="00E27454"><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:bookmarkStart w:id="4" w:name="_Toc492960772"/><w:r w:rsidRPr="00E27454"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t>Interfaces</w:t></w:r><w:bookmarkEnd w:id="4"/><w:r w:rsidRPr="00E27454"><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> to External Systems or Devices</w:t></w:r></w:p>'
        $p.Range.InsertXML($pkgHead + $xml + $pkgTail)
        Write-Host "  done."
        break
    }
}

$d.Save()
Write-Host "Saved."
